$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between the two detail rows (16 and 17).
$ws.Range("E16").Value = "2001"
$ws.Range("E17").Value = "2002"

# Update the "Valor Mora" amounts for both detail rows.
$ws.Range("G16").Value = 1350000
$ws.Range("G17").Value = 1350000
